$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 15.140316
$ws.Range("H2").Value = 45.420948
$ws.Range("I2").Value = 0.3229157245229468
$ws.Range("J2").Value = 0.3229157245229468
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 18.76192366666666
$ws.Range("N2").Value = 56.285771
$ws.Range("O2").Value = 0.1222461152048115
$ws.Range("P2").Value = 0.1222461152048115
$ws.Range("Q2").Value = 284.0614530812119
$ws.Range("R2").Value = 2556.553077730908
$ws.Range("S2").Value = 0.03947519286147733
$ws.Range("T2").Value = 0.03947519286147733
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 15.140316
$ws.Range("H3").Value = 45.420948
$ws.Range("I3").Value = 0.3229157245229468
$ws.Range("J3").Value = 0.3229157245229468
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 38.59812166666666
$ws.Range("N3").Value = 115.794365
$ws.Range("O3").Value = 0.2514918252404857
$ws.Range("P3").Value = 0.2514918252404857
$ws.Range("Q3").Value = 584.3877590397799
$ws.Range("R3").Value = 5259.48983135802
$ws.Range("S3").Value = 0.08121066495912975
$ws.Range("T3").Value = 0.08121066495912975
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 15.140316
$ws.Range("H4").Value = 45.420948
$ws.Range("I4").Value = 0.3229157245229468
$ws.Range("J4").Value = 0.3229157245229468
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.691683333333334
$ws.Range("N4").Value = 20.07505
$ws.Range("O4").Value = 0.04360066196912097
$ws.Range("P4").Value = 0.04360066196912097
$ws.Range("Q4").Value = 101.3142002386
$ws.Range("R4").Value = 911.8278021474
$ws.Range("S4").Value = 0.01407933934943879
$ws.Range("T4").Value = 0.01407933934943879
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 15.140316
$ws.Range("H5").Value = 45.420948
$ws.Range("I5").Value = 0.3229157245229468
$ws.Range("J5").Value = 0.3229157245229468
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 89.42491666666666
$ws.Range("N5").Value = 268.27475
$ws.Range("O5").Value = 0.5826613975855819
$ws.Range("P5").Value = 0.5826613975855818
$ws.Range("Q5").Value = 1353.921496607
$ws.Range("R5").Value = 12185.293469463
$ws.Range("S5").Value = 0.1881505273529009
$ws.Range("T5").Value = 0.1881505273529009
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 18.94069966666667
$ws.Range("H6").Value = 56.822099
$ws.Range("I6").Value = 0.4039710766824948
$ws.Range("J6").Value = 0.4039710766824948
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 18.76192366666666
$ws.Range("N6").Value = 56.285771
$ws.Range("O6").Value = 0.1222461152048115
$ws.Range("P6").Value = 0.1222461152048115
$ws.Range("Q6").Value = 355.3639613392588
$ws.Range("R6").Value = 3198.275652053329
$ws.Range("S6").Value = 0.04938389477954001
$ws.Range("T6").Value = 0.04938389477954001
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 18.94069966666667
$ws.Range("H7").Value = 56.822099
$ws.Range("I7").Value = 0.4039710766824948
$ws.Range("J7").Value = 0.4039710766824948
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 38.59812166666666
$ws.Range("N7").Value = 115.794365
$ws.Range("O7").Value = 0.2514918252404857
$ws.Range("P7").Value = 0.2514918252404857
$ws.Range("Q7").Value = 731.0754301857927
$ws.Range("R7").Value = 6579.678871672135
$ws.Range("S7").Value = 0.1015954234192448
$ws.Range("T7").Value = 0.1015954234192448
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 18.94069966666667
$ws.Range("H8").Value = 56.822099
$ws.Range("I8").Value = 0.4039710766824948
$ws.Range("J8").Value = 0.4039710766824948
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.691683333333334
$ws.Range("N8").Value = 20.07505
$ws.Range("O8").Value = 0.04360066196912097
$ws.Range("P8").Value = 0.04360066196912097
$ws.Range("Q8").Value = 126.7451642811056
$ws.Range("R8").Value = 1140.70647852995
$ws.Range("S8").Value = 0.0176134063597353
$ws.Range("T8").Value = 0.0176134063597353
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 18.94069966666667
$ws.Range("H9").Value = 56.822099
$ws.Range("I9").Value = 0.4039710766824948
$ws.Range("J9").Value = 0.4039710766824948
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 89.42491666666666
$ws.Range("N9").Value = 268.27475
$ws.Range("O9").Value = 0.5826613975855819
$ws.Range("P9").Value = 0.5826613975855818
$ws.Range("Q9").Value = 1693.770489300028
$ws.Range("R9").Value = 15243.93440370025
$ws.Range("S9").Value = 0.2353783521239747
$ws.Range("T9").Value = 0.2353783521239746
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.221232
$ws.Range("H10").Value = 9.663696
$ws.Range("I10").Value = 0.06870308817441464
$ws.Range("J10").Value = 0.06870308817441464
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 18.76192366666666
$ws.Range("N10").Value = 56.285771
$ws.Range("O10").Value = 0.1222461152048115
$ws.Range("P10").Value = 0.1222461152048115
$ws.Range("Q10").Value = 60.436508896624
$ws.Range("R10").Value = 543.9285800696159
$ws.Range("S10").Value = 0.008398685631895818
$ws.Range("T10").Value = 0.008398685631895818
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 3.221232
$ws.Range("H11").Value = 9.663696
$ws.Range("I11").Value = 0.06870308817441464
$ws.Range("J11").Value = 0.06870308817441464
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 38.59812166666666
$ws.Range("N11").Value = 115.794365
$ws.Range("O11").Value = 0.2514918252404857
$ws.Range("P11").Value = 0.2514918252404857
$ws.Range("Q11").Value = 124.33350465256
$ws.Range("R11").Value = 1119.00154187304
$ws.Range("S11").Value = 0.01727826504464157
$ws.Range("T11").Value = 0.01727826504464157
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 3.221232
$ws.Range("H12").Value = 9.663696
$ws.Range("I12").Value = 0.06870308817441464
$ws.Range("J12").Value = 0.06870308817441464
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 6.691683333333334
$ws.Range("N12").Value = 20.07505
$ws.Range("O12").Value = 0.04360066196912097
$ws.Range("P12").Value = 0.04360066196912097
$ws.Range("Q12").Value = 21.5554644872
$ws.Range("R12").Value = 193.9991803848
$ws.Range("S12").Value = 0.002995500123727365
$ws.Range("T12").Value = 0.002995500123727365
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 3.221232
$ws.Range("H13").Value = 9.663696
$ws.Range("I13").Value = 0.06870308817441464
$ws.Range("J13").Value = 0.06870308817441464
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 89.42491666666666
$ws.Range("N13").Value = 268.27475
$ws.Range("O13").Value = 0.5826613975855819
$ws.Range("P13").Value = 0.5826613975855818
$ws.Range("Q13").Value = 288.058403164
$ws.Range("R13").Value = 2592.525628476
$ws.Range("S13").Value = 0.0400306373741499
$ws.Range("T13").Value = 0.04003063737414989
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 9.584029000000001
$ws.Range("H14").Value = 28.752087
$ws.Range("I14").Value = 0.2044101106201438
$ws.Range("J14").Value = 0.2044101106201438
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 18.76192366666666
$ws.Range("N14").Value = 56.285771
$ws.Range("O14").Value = 0.1222461152048115
$ws.Range("P14").Value = 0.1222461152048115
$ws.Range("Q14").Value = 179.8148205171197
$ws.Range("R14").Value = 1618.333384654077
$ws.Range("S14").Value = 0.02498834193189837
$ws.Range("T14").Value = 0.02498834193189837
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 9.584029000000001
$ws.Range("H15").Value = 28.752087
$ws.Range("I15").Value = 0.2044101106201438
$ws.Range("J15").Value = 0.2044101106201438
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 38.59812166666666
$ws.Range("N15").Value = 115.794365
$ws.Range("O15").Value = 0.2514918252404857
$ws.Range("P15").Value = 0.2514918252404857
$ws.Range("Q15").Value = 369.9255173988617
$ws.Range("R15").Value = 3329.329656589755
$ws.Range("S15").Value = 0.05140747181746955
$ws.Range("T15").Value = 0.05140747181746956
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 9.584029000000001
$ws.Range("H16").Value = 28.752087
$ws.Range("I16").Value = 0.2044101106201438
$ws.Range("J16").Value = 0.2044101106201438
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 6.691683333333334
$ws.Range("N16").Value = 20.07505
$ws.Range("O16").Value = 0.04360066196912097
$ws.Range("P16").Value = 0.04360066196912097
$ws.Range("Q16").Value = 64.13328712548335
$ws.Range("R16").Value = 577.1995841293501
$ws.Range("S16").Value = 0.008912416136219514
$ws.Range("T16").Value = 0.008912416136219516
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 9.584029000000001
$ws.Range("H17").Value = 28.752087
$ws.Range("I17").Value = 0.2044101106201438
$ws.Range("J17").Value = 0.2044101106201438
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 89.42491666666666
$ws.Range("N17").Value = 268.27475
$ws.Range("O17").Value = 0.5826613975855819
$ws.Range("P17").Value = 0.5826613975855818
$ws.Range("Q17").Value = 857.0509946559167
$ws.Range("R17").Value = 7713.458951903251
$ws.Range("S17").Value = 0.1191018807345564
$ws.Range("T17").Value = 0.1191018807345564

Write-Output "Done: updated 224 cells"
